$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newK = @{
    2  = 1
    3  = 3
    5  = 2
    6  = 0
    7  = 3
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 0
    15 = 3
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 2
    21 = 3
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 3
    31 = 3
    32 = 0
    33 = 0
    34 = 3
    35 = 1
    36 = 3
    37 = 1
    38 = 6
    39 = 0
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 3
    45 = 3
    46 = 2
    47 = 2
    48 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
